# Updated result slides
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename "Comparison - AUPRC (2)" -> "Comparison - Ablation"
# ---------------------------------------------------------------------
$wsAblation = $wb.Worksheets.Item("Comparison - AUPRC (2)")
$wsAblation.Name = "Comparison - Ablation"

# ---------------------------------------------------------------------
# 2. Fix the "Comparision %" -> "Comparison %" typo everywhere it
#    occurs so the shared-string table collapses/compacts exactly like
#    it does in the authored workbook.
# ---------------------------------------------------------------------
$sheetNames = @(
    "Comparison - AUPRC",
    "Comparison - AUROC",
    "Comparison - vs multi",
    "Comparison - vs timeseries",
    "Comparison - vs ablation",
    "Comparison - Ablation"
)
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $rng = $ws.UsedRange
    for ($r = 1; $r -le $rng.Rows.Count; $r++) {
        for ($c = 1; $c -le $rng.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value2 -eq "Comparision %") {
                $cell.Value = "Comparison %"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 3. Relabel the two comparison rows in "Comparison - Ablation" so
#    they compare against the Best Baseline instead of Multimodal.
# ---------------------------------------------------------------------
$wsAblation.Range("A6").Value = "Proposed vs Best Baseline"
$wsAblation.Range("A7").Value = "Glove vs Best Baseline"

# Column A widens (still best-fit) to accommodate the new, longer labels.
$wsAblation.Columns.Item(1).ColumnWidth = 21.998697916666668

# ---------------------------------------------------------------------
# 4. Restore each sheet's on-screen selection to match the saved view.
# ---------------------------------------------------------------------
$wsAuroc = $wb.Worksheets.Item("Comparison - AUROC")
$wsAuroc.Activate()
$wsAuroc.Range("J31").Select() | Out-Null

$wsMulti = $wb.Worksheets.Item("Comparison - vs multi")
$wsMulti.Activate()
$wsMulti.Range("G10").Select() | Out-Null

$wsTime = $wb.Worksheets.Item("Comparison - vs timeseries")
$wsTime.Activate()
$wsTime.Range("A1:E5").Select() | Out-Null

$wsVsAblation = $wb.Worksheets.Item("Comparison - vs ablation")
$wsVsAblation.Activate()
$wsVsAblation.Range("A5").Select() | Out-Null

# Leave "Comparison - Ablation" as the active / selected tab, matching
# the workbook's saved activeTab and the sheet's own tabSelected flag.
$wsAblation.Activate()
$wsAblation.Range("A8:E8").Select() | Out-Null
